{"js": "// The document contained a Quarto/Pandoc conditional-publishing div:\n//   ::: {.content-visble when-format=\"pdf\", \"pptx\", \"docx\"}\n//   Links to Viktoria's portfolio pieces are only available in the html output.\n//   :::\n// split across many runs/paragraphs. The edit removes the div \"fence\"\n// paragraphs (the opening \"::: {...}\" line and the closing \":::\" line)\n// and promotes the remaining sentence back onto the \"FirstParagraph\"\n// style that the opening fence paragraph used to carry.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the Pandoc fenced-div delimiter paragraphs (\":::\" at the start\n// of the opening line, and the whole closing line).\nconst fenceParagraphs = paragraphs.items.filter((p) =>\n  p.text.trim().startsWith(\":::\")\n);\n\nfor (const p of fenceParagraphs) {\n  p.delete();\n}\nawait context.sync();\n\n// The sentence that used to sit between the fences becomes its own\n// paragraph again; restore the \"FirstParagraph\" style it should carry.\nconst remaining = body.paragraphs;\nremaining.load(\"items/text\");\nawait context.sync();\n\nconst linksParagraph = remaining.items.find((p) =>\n  p.text.indexOf(\"Links to Viktoria\") !== -1\n);\nif (linksParagraph) {\n  linksParagraph.style = \"FirstParagraph\";\n}\nawait context.sync();\n", "ps1": "# The document contained a Quarto/Pandoc conditional-publishing div:\n#   ::: {.content-visble when-format=\"pdf\", \"pptx\", \"docx\"}\n#   Links to Viktoria's portfolio pieces are only available in the html output.\n#   :::\n# The edit removes the div \"fence\" paragraphs (the opening \"::: {...}\"\n# line and the closing \":::\" line) and promotes the remaining sentence\n# back onto the \"FirstParagraph\" style that the opening fence paragraph\n# used to carry.\n\n$d = $word.ActiveDocument\n\n# Find the Pandoc fenced-div delimiter paragraphs (\":::\" at the start\n# of the opening line, and the whole closing line).\n$toDelete = @()\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.Trim()\n    if ($t.StartsWith(\":::\")) {\n        $toDelete += $p\n    }\n}\n\n# Delete from the last match to the first so earlier Range objects in\n# the list don't get invalidated by removing a paragraph that follows\n# them in the document.\nfor ($i = $toDelete.Count - 1; $i -ge 0; $i--) {\n    $toDelete[$i].Range.Delete()\n}\n\n# The sentence that used to sit between the fences becomes its own\n# paragraph again; restore the \"FirstParagraph\" style it should carry.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Links to Viktoria*\") {\n        $p.Style = \"FirstParagraph\"\n    }\n}\n"}
